$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = -0.005020920502092054
$ws.Range("H2").Value = -0.06631799163179916
$ws.Range("I2").Value = -0.1077405857740586
$ws.Range("J2").Value = -0.1077405857740586
$ws.Range("K2").Value = -1.17
$ws.Range("L2").Value = -0.1223849372384937
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 0.105
$ws.Range("V2").Value = 0.01761744966442953
$ws.Range("W2").Value = -0.740506329113924
$ws.Range("X2").Value = 0.08494330876030783
$ws.Range("Y2").Value = -0.8254496378742319
$ws.Range("Z2").Value = 3.64607170099161
$ws.Range("AA2").Value = -0.3928299008390542
$ws.Range("AB2").Value = 0.07264391290481806
$ws.Range("AC2").Value = -0.4654738137438722
$ws.Range("AD2").Value = 1.49
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 1.49
$ws.Range("AG2").Value = 1.385
$ws.Range("AH2").Value = 0.2
$ws.Range("AI2").Value = 0.6263135771332493
$ws.Range("AJ2").Value = 0.1885636487406399
$ws.Range("AK2").Value = 0.6090589270008795
$ws.Range("AL2").Value = 0.152
$ws.Range("AM2").Value = 0.151
$ws.Range("AN2").Value = -2.704174228675136
$ws.Range("AO2").Value = -6.776315789473685
$ws.Range("AP2").Value = -2.513611615245009
$ws.Range("AQ2").Value = -6.821192052980133

# Row 3 updates
$ws.Range("G3").Value = -0.005020920502092054
$ws.Range("H3").Value = -0.06631799163179916
$ws.Range("I3").Value = -0.1077405857740586
$ws.Range("J3").Value = -0.1077405857740586
$ws.Range("K3").Value = -1.17
$ws.Range("L3").Value = -0.1223849372384937
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 0.105
$ws.Range("V3").Value = 0.01761744966442953
$ws.Range("W3").Value = -0.740506329113924
$ws.Range("X3").Value = 0.08494330876030783
$ws.Range("Y3").Value = -0.8254496378742319
$ws.Range("Z3").Value = 3.64607170099161
$ws.Range("AA3").Value = -0.3928299008390542
$ws.Range("AB3").Value = 0.07264391290481806
$ws.Range("AC3").Value = -0.4654738137438722
$ws.Range("AD3").Value = 1.49
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 1.49
$ws.Range("AG3").Value = 1.385
$ws.Range("AH3").Value = 0.2
$ws.Range("AI3").Value = 0.6263135771332493
$ws.Range("AJ3").Value = 0.1885636487406399
$ws.Range("AK3").Value = 0.6090589270008795
$ws.Range("AL3").Value = 0.152
$ws.Range("AM3").Value = 0.151
$ws.Range("AN3").Value = -2.704174228675136
$ws.Range("AO3").Value = -6.776315789473685
$ws.Range("AP3").Value = -2.513611615245009
$ws.Range("AQ3").Value = -6.821192052980133
